$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value  = 16.1669
$ws.Range("D9").Value  = -7.379400000000001
$ws.Range("E12").Value = 18.06590000000003
$ws.Range("E14").Value = 16.80510000000001
$ws.Range("D18").Value = -8.576499999999992
$ws.Range("D20").Value = -7.420899999999998
$ws.Range("E26").Value = 16.2256
$ws.Range("D27").Value = -8.385499999999997
$ws.Range("E27").Value = 16.77289999999999
$ws.Range("E29").Value = 16.95280000000001
$ws.Range("D35").Value = -7.891699999999999
$ws.Range("E37").Value = 16.65760000000001
$ws.Range("E38").Value = 16.63719999999999
$ws.Range("E51").Value = 17.2412
$ws.Range("E52").Value = 16.8604
$ws.Range("E55").Value = 16.4071
$ws.Range("D69").Value = -7.248799999999997
$ws.Range("E69").Value = 17.16090000000002
$ws.Range("E70").Value = 18.03670000000002
$ws.Range("D76").Value = -7.690599999999998
$ws.Range("D78").Value = -7.690500000000006
$ws.Range("E81").Value = 16.55579999999999
$ws.Range("D82").Value = -8.317199999999991
$ws.Range("D83").Value = -8.953300000000002
$ws.Range("E83").Value = 16.68869999999998
$ws.Range("D93").Value = -6.351199999999997
$ws.Range("E102").Value = 16.8931
